# Apply cryptos list update (prices & 1h volume deltas) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, even when it looks like a number
# (e.g. "1.00", "597.58"), so Excel does not silently convert it to a
# numeric cell. Forcing the Text number format before assignment, then
# clearing formats afterwards, keeps the cell style untouched.
function Set-TextValue($cellAddr, $text) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

$ws.Range("D2").Value = "63.758.16"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.620.72"
$ws.Range("E3").Value = "  -0.08%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "597.58"
$ws.Range("E5").Value = "  -1.50%  "
Set-TextValue "D6" "151.07"
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("E7").Value = "  +0.12%  "
Set-TextValue "D8" "0.591"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("E9").Value = "  +0.49%  "
Set-TextValue "D10" "5.70"
$ws.Range("E10").Value = "  +2.74%  "
Set-TextValue "D11" "0.385"
$ws.Range("E11").Value = "  +3.49%  "
$ws.Range("E12").Value = "  -0.96%  "
Set-TextValue "D13" "27.87"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "3.091.86"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "63.563.82"
$ws.Range("E15").Value = "  +0.52%  "
Set-TextValue "D16" "0.0000153"
$ws.Range("E16").Value = "  +4.21%  "
$ws.Range("D17").Value = "2.622.10"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("E18").Value = "  +6.94%  "
Set-TextValue "D19" "4.70"
$ws.Range("E19").Value = "  +3.25%  "
Set-TextValue "D20" "347.81"
$ws.Range("E20").Value = "  +1.41%  "
Set-TextValue "D21" "6.88"
$ws.Range("E21").Value = "  +0.30%  "
Set-TextValue "D22" "0.999"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  +1.33%  "
Set-TextValue "D24" "66.81"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  +8.48%  "
Set-TextValue "D26" "9.32"
$ws.Range("E26").Value = "  +2.48%  "
Set-TextValue "D27" "1.68"
$ws.Range("E27").Value = "  -1.45%  "
Set-TextValue "D28" "555.52"
$ws.Range("E28").Value = "  -0.66%  "
Set-TextValue "D29" "8.16"
$ws.Range("E29").Value = "  +3.72%  "
Set-TextValue "D30" "0.162"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  -0.07%  "
Set-TextValue "D32" "2.05"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "0.0₃0851"
$ws.Range("E33").Value = "  +0.07%  "
Set-TextValue "D34" "1.77"
$ws.Range("E34").Value = "  +1.10%  "
Set-TextValue "D35" "5.28"
$ws.Range("E35").Value = "  +2.45%  "
Set-TextValue "D36" "168.10"
$ws.Range("E36").Value = "  +0.54%  "
Set-TextValue "D37" "0.415"
$ws.Range("E37").Value = "  +2.99%  "
$ws.Range("E38").Value = "  -0.26%  "
Set-TextValue "D39" "19.51"
$ws.Range("E39").Value = "  +2.47%  "
Set-TextValue "D40" "1.94"
$ws.Range("E40").Value = "  +0.26%  "
Set-TextValue "D41" "0.999"
$ws.Range("E41").Value = "  -0.04%  "
Set-TextValue "D42" "166.60"
$ws.Range("E42").Value = "  +0.68%  "
Set-TextValue "D43" "39.73"
$ws.Range("E43").Value = "  -0.11%  "
Set-TextValue "D44" "3.94"
$ws.Range("E44").Value = "  +4.11%  "
Set-TextValue "D45" "0.0591"
$ws.Range("E45").Value = "  +4.12%  "
Set-TextValue "D46" "21.71"
$ws.Range("E46").Value = "  -0.54%  "
Set-TextValue "D47" "0.629"
$ws.Range("E47").Value = "  +0.31%  "
Set-TextValue "D48" "0.0252"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D49" "2.01"
$ws.Range("E49").Value = "  +4.42%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0249"
$ws.Range("E50").Value = "  +25.88%  "
Set-TextValue "D51" "0.0967"
$ws.Range("E51").Value = "  +0.87%  "
